$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mini 3x3 neighbor grid around the "*" cell (F15)
$ws.Range("E14").Value = "(r-1, c-1)"
$ws.Range("F14").Value = "(r-1, c)"
$ws.Range("G14").Value = "(r-1, c+1)"
$ws.Range("E15").Value = "(r, c-1)"
$ws.Range("F15").Value = "* (r, c)"
$ws.Range("G15").Value = "(r, c+1)"
$ws.Range("E16").Value = "(r+1, c-1)"
$ws.Range("F16").Value = "(r+1, c)"
$ws.Range("G16").Value = "(r+1, c+1)"

# Numbered list (I column) matched with neighbor coordinate labels (J column)
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = "(r-1, c-1)"
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = "(r-1, c)"
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = "(r-1, c+1)"
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = "(r, c-1)"
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = "(r, c+1)"
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = "(r+1, c-1)"
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = "(r+1, c)"
$ws.Range("I21").Value = 8
$ws.Range("J21").Value = "(r+1, c+1)"

# Heading for the neighbor list
$ws.Range("J13").Value = "Calculate potnetial neighbors of (r, c)"

# Rules notes (N column)
$ws.Range("N13").Value = "Rules to validate neighbor points"
$ws.Range("N14").Value = "1) r of neighbor must be <= row len"
$ws.Range("N15").Value = "2) c of neighbor must be <= column len"
$ws.Range("N16").Value = "3) other mines don't count as neighbors"

# Restore the selection to match the saved workbook state
$ws.Range("O28").Select()
